# Atualização automática da planilha
# Target sheet: "Organograma"
#  - Row 6 (Comitê Executivo / Membro do Comitê) person changed from
#    "Paulo Figueiredo" (Diretor Financeiro / Financeiro) to
#    "Henrique Hildebrand Garcia" (Diretor Jurídico / Jurídico).
#  - A new blank row is inserted right after (new row 7), pushing every
#    subsequent row down by one (old row 7 -> new row 8, ..., old row 27 ->
#    new row 28), which also grows the used range / AutoFilter / hidden
#    filter-database defined name from A1:E27 to A1:E28.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Organograma")

# 1) Update the three changed cells on row 6.
$ws.Range("C6").Value = "Henrique Hildebrand Garcia"
$ws.Range("D6").Value = "Diretor Jurídico"
$ws.Range("E6").Value = "Jurídico"
$ws.Rows.Item(6).RowHeight = 15

# 2) Insert a brand-new blank row directly below row 6 (becomes row 7),
#    shifting rows 7..27 down to 8..28. Give it the same formatting as the
#    row above it, but leave its cells empty.
$ws.Rows.Item(7).Insert()
$ws.Range("A6:E6").Copy()
$ws.Range("A7:E7").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Rows.Item(7).RowHeight = 15

# 3) The sheet's AutoFilter needs to grow from A1:E27 to A1:E28 to include
#    the newly added row. Toggle off/on so the new extent is picked up.
$ws.AutoFilterMode = $false
$ws.Range("A1:E28").AutoFilter()

# 4) The hidden _FilterDatabase defined name also tracks the filtered
#    range and must be expanded to match.
$fdb = $wb.Names.Item("Organograma!_FilterDatabase")
$fdb.RefersTo = "=Organograma!`$A`$1:`$E`$28"

# 5) Leave the sheet's selection on the edited row, matching the saved
#    cursor position of the authored change.
$ws.Rows.Item(6).Select()
